$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1969163036.699543
$ws.Range("C3").Value = 2602886220.450687
$ws.Range("B4").Value = 5538980.330819745
$ws.Range("C4").Value = 7100049.906343766
